# Guion 02 grado tercero - corrección de estilo
# Applies the textual corrections described in the commit:
#  - "Aplicar" -> "Aplica" (x2)
#  - "Sistema de numeración decimal,lectura," -> adds missing spaces
#  - " adultos mientras..." -> adds missing comma
#  - "...docente vía correo..." -> "...docente por correo electrónico..." (x2)
#  - "ubicada" -> "ubicado" (gender agreement) + relocate the _GoBack bookmark
#  - merges a couple of runs that only differed by an accidental split

$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1) "Aplicar el sistema..." -> "Aplica el sistema..." (both occurrences in the doc)
Replace-All "Aplicar el sistema de numeración decimal en diferentes contextos" `
            "Aplica el sistema de numeración decimal en diferentes contextos"

# 2) Missing spaces in the keyword list
Replace-All "Sistema de numeración decimal,lectura," `
            "Sistema de numeración decimal, lectura, "

# 3) Merge the accidentally split filename runs (content itself is unchanged)
Replace-All "MA_03_02_CO_REC100_IMG01n" "MA_03_02_CO_REC100_IMG01n"
Replace-All "MA_03_02_CO_REC100_IMG02n" "MA_03_02_CO_REC100_IMG02n"

# 4) Missing comma
Replace-All " adultos mientras que en el mes de mayo ingresaron " `
            " adultos, mientras que en el mes de mayo ingresaron "

# 5) "vía correo" -> "por correo electrónico" (first occurrence, single run, extra
#    trailing spaces cleaned up)
Replace-All "Escribe en letras las cantidades que se nombran en la situación y envía tu respuesta al docente vía correo o entrégala en una hoja.  " `
            "Escribe en letras las cantidades que se nombran en la situación y envía tu respuesta al docente por correo electrónico o entrégala en una hoja."

# 6) Same sentence fix, second occurrence (spans two runs, no trailing spaces here)
Replace-All "Escribe en letras las cantidades que se nombran en la situación y envía tu respuesta al docente vía correo o entrégala en una hoja." `
            "Escribe en letras las cantidades que se nombran en la situación y envía tu respuesta al docente por correo electrónico o entrégala en una hoja."

# 7) Gender agreement fix: "ubicada" -> "ubicado"
Replace-All "Ojos del Salado ubicada en Chile con" "Ojos del Salado ubicado en Chile con"

# 8) Relocate the "_GoBack" bookmark to just after "ubicado" (where Word last left the
#    cursor), removing it from its old spot after "...en una hoja."
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Ojos del Salado ubicado"
$found = $rng.Find.Execute()
if ($found) {
    $bmPos = $rng.End
    $bmRng = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRng)
}

Write-Output "done"
